# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (want-to-go count) figures and marks a
# now-unavailable ticket price as "不可售" across the relevant sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet1) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 281
$ws.Range("G3").Value = "不可售"
$ws.Range("F6").Value = 658
$ws.Range("F12").Value = 3362
$ws.Range("F14").Value = 75
$ws.Range("F16").Value = 36
$ws.Range("F17").Value = 50
$ws.Range("F18").Value = 571
$ws.Range("F19").Value = 42
$ws.Range("F20").Value = 668
$ws.Range("F26").Value = 2382
$ws.Range("F27").Value = 4925
$ws.Range("F31").Value = 1263
$ws.Range("F32").Value = 268
$ws.Range("F33").Value = 2187
$ws.Range("F35").Value = 483
$ws.Range("F41").Value = 770

# ---- Sheet "演出" (sheet2) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 65

# ---- Sheet "全部类型" (sheet4) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 281
$ws.Range("G3").Value = "不可售"
$ws.Range("F6").Value = 658
$ws.Range("F12").Value = 3362
$ws.Range("F14").Value = 75
$ws.Range("F16").Value = 65
$ws.Range("F17").Value = 36
$ws.Range("F18").Value = 50
$ws.Range("F19").Value = 571
$ws.Range("F20").Value = 42
$ws.Range("F21").Value = 668
$ws.Range("F27").Value = 2382
$ws.Range("F28").Value = 4925
$ws.Range("F32").Value = 1263
$ws.Range("F33").Value = 268
$ws.Range("F34").Value = 2187
$ws.Range("F36").Value = 483
$ws.Range("F42").Value = 770
